# The deck currently uses the "Integral" (Red Violet) theme colours on its
# one-and-only slide master/design. Re-point every theme colour to the
# stock "Office Theme" palette (the colours that were sitting, unused, in
# the deck's spare theme part) so the design's colour scheme becomes the
# plain Office one.

$p  = $ppt.ActivePresentation
$m  = $p.SlideMaster
$th = $m.Theme
$cs = $th.ThemeColorScheme

function Set-ThemeColor($index, $hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $cs.Item($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Index order matches <a:clrScheme>: dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink
Set-ThemeColor 1  "000000"
Set-ThemeColor 2  "FFFFFF"
Set-ThemeColor 3  "44546A"
Set-ThemeColor 4  "E7E6E6"
Set-ThemeColor 5  "5B9BD5"
Set-ThemeColor 6  "ED7D31"
Set-ThemeColor 7  "A5A5A5"
Set-ThemeColor 8  "FFC000"
Set-ThemeColor 9  "4472C4"
Set-ThemeColor 10 "70AD47"
Set-ThemeColor 11 "0563C1"
Set-ThemeColor 12 "954F72"
